$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# The header row gains a "Recipient" and "myusername" column (inserted
# among the existing fields) and a new trailing "Email Sent" column.
# The previous "myemail" column is removed entirely.
#
# Resulting layout (A1:I1):
#   A1 fname | B1 lname | C1 Recipient | D1 myphone | E1 myusername |
#   F1 mypassword | G1 mybsb | H1 myaccount | I1 Email Sent
# ---------------------------------------------------------------------

$ws.Range("A1").Value = "fname"
$ws.Range("B1").Value = "lname"
$ws.Range("C1").Value = "Recipient"
$ws.Range("D1").Value = "myphone"
$ws.Range("E1").Value = "myusername"
$ws.Range("F1").Value = "mypassword"
$ws.Range("G1").Value = "mybsb"
$ws.Range("H1").Value = "myaccount"

# I1 is a brand new cell; give it the same cell style the rest of the
# header row already uses (s="1") before writing its text, by copying
# the formatting from the neighbouring H1 cell.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "Email Sent"
